$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 390; this shifts existing rows 390..440 down to 391..441
$ws.Rows.Item(390).Insert()

# Fill the new row 390 with its data
$ws.Cells.Item(390, 1).Value = 5
$ws.Cells.Item(390, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(390, 3).Value = "Maule"
$ws.Cells.Item(390, 4).Value = 44984
$ws.Cells.Item(390, 5).Value = 7
$ws.Cells.Item(390, 6).Value = 100112003
$ws.Cells.Item(390, 7).Value = "Ajo"
$ws.Cells.Item(390, 8).Value = "Chino"
$ws.Cells.Item(390, 9).Value = "Primera"
$ws.Cells.Item(390, 10).Value = 200
$ws.Cells.Item(390, 11).Value = 20000
$ws.Cells.Item(390, 12).Value = 20000
$ws.Cells.Item(390, 13).Value = 20000
$ws.Cells.Item(390, 14).Value = "$/malla 10 kilos"
$ws.Cells.Item(390, 15).Value = "China"
$ws.Cells.Item(390, 16).Value = 2000
$ws.Cells.Item(390, 17).Value = 10
$ws.Cells.Item(390, 18).Value = "Hortaliza"
